$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated crypto price (column D) and volume-change (column E) values,
# one entry per data row (rows 2-51). Where D is $null, the price value
# is unchanged and only the volume percentage is updated.
$updates = @(
    @{ Row = 2; D = "27.410.87"; E = "  -2.69%  " }
    @{ Row = 3; D = "1.743.27"; E = "  -3.28%  " }
    @{ Row = 4; D = $null; E = "  -0.08%  " }
    @{ Row = 5; D = "321.24"; E = "  -4.46%  " }
    @{ Row = 6; D = $null; E = "  +0.01%  " }
    @{ Row = 7; D = "0.4200"; E = "  -9.20%  " }
    @{ Row = 8; D = "0.3581"; E = "  -3.09%  " }
    @{ Row = 9; D = "45.48"; E = "  +0.58%  " }
    @{ Row = 10; D = "0.07409"; E = "  -2.60%  " }
    @{ Row = 11; D = "1.113"; E = "  -3.08%  " }
    @{ Row = 12; D = $null; E = "  -0.07%  " }
    @{ Row = 13; D = $null; E = "  -4.04%  " }
    @{ Row = 14; D = "6.092"; E = "  -4.00%  " }
    @{ Row = 15; D = "7.184"; E = "  -3.59%  " }
    @{ Row = 16; D = "1.738.23"; E = "  -3.49%  " }
    @{ Row = 17; D = $null; E = "  -2.86%  " }
    @{ Row = 18; D = "87.52"; E = "  +6.85%  " }
    @{ Row = 19; D = "0.06051"; E = "  -9.90%  " }
    @{ Row = 20; D = $null; E = "  +0.02%  " }
    @{ Row = 21; D = "16.85"; E = "  -3.79%  " }
    @{ Row = 22; D = "6.102"; E = "  -4.94%  " }
    @{ Row = 23; D = "0.5233"; E = "  -5.47%  " }
    @{ Row = 24; D = "27.428.25"; E = "  -2.62%  " }
    @{ Row = 25; D = "11.41"; E = "  -4.19%  " }
    @{ Row = 26; D = "2.339"; E = "  -2.98%  " }
    @{ Row = 27; D = "20.41"; E = "  -1.87%  " }
    @{ Row = 28; D = "152.91"; E = "  -0.30%  " }
    @{ Row = 29; D = "2.379"; E = "  -0.06%  " }
    @{ Row = 30; D = "1.936.32"; E = "  -3.49%  " }
    @{ Row = 31; D = "125.66"; E = "  -5.49%  " }
    @{ Row = 32; D = "1.180"; E = "  -5.95%  " }
    @{ Row = 33; D = "5.684"; E = "  -2.99%  " }
    @{ Row = 34; D = "0.09112"; E = "  -4.47%  " }
    @{ Row = 35; D = "3.633"; E = "  -9.85%  " }
    @{ Row = 36; D = "12.66"; E = "  +4.83%  " }
    @{ Row = 37; D = "0.02290"; E = "  -2.70%  " }
    @{ Row = 38; D = "0.2133"; E = "  -4.31%  " }
    @{ Row = 39; D = "5.075"; E = "  -3.38%  " }
    @{ Row = 40; D = "0.06051"; E = "  -5.15%  " }
    @{ Row = 41; D = "0.6381"; E = "  -4.06%  " }
    @{ Row = 42; D = "1.193"; E = "  -3.98%  " }
    @{ Row = 43; D = "1.419"; E = "  -6.22%  " }
    @{ Row = 44; D = "1.0000"; E = "  +0.04%  " }
    @{ Row = 45; D = "7.924"; E = "  -2.19%  " }
    @{ Row = 46; D = "13.68"; E = "  -3.04%  " }
    @{ Row = 47; D = "3.704"; E = "  -3.37%  " }
    @{ Row = 48; D = "0.5830"; E = "  -4.58%  " }
    @{ Row = 49; D = "125.32"; E = "  -3.84%  " }
    @{ Row = 50; D = "1.945"; E = "  -5.34%  " }
    @{ Row = 51; D = "0.06821"; E = "  -4.56%  " }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        $ws.Cells.Item($u.Row, 4).Value = $u.D
    }
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}
